$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 311-312; existing rows 311-334 shift down to 313-336.
$ws.Rows("311:312").Insert()

# New row 311 - Apio, Primera, fecha 44610
$ws.Range("A311").Value = 8
$ws.Range("B311").Value = "Terminal La Palmera de La Serena"
$ws.Range("C311").Value = "Coquimbo"
$ws.Range("D311").Value = 44610
$ws.Range("D311").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E311").Value = 4
$ws.Range("F311").Value = 100112017
$ws.Range("G311").Value = "Apio"
$ws.Range("H311").Value = "Americana (o)"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 2000
$ws.Range("K311").Value = 7500
$ws.Range("L311").Value = 8000
$ws.Range("M311").Value = 7750
$ws.Range("N311").Value = "`$/docena de matas"
$ws.Range("O311").Value = "Provincia del Elquí"
$ws.Range("P311").Value = 1292
$ws.Range("Q311").Value = 6
$ws.Range("R311").Value = "Hortaliza"

# New row 312 - Apio, Segunda, fecha 44610
$ws.Range("A312").Value = 8
$ws.Range("B312").Value = "Terminal La Palmera de La Serena"
$ws.Range("C312").Value = "Coquimbo"
$ws.Range("D312").Value = 44610
$ws.Range("D312").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E312").Value = 4
$ws.Range("F312").Value = 100112017
$ws.Range("G312").Value = "Apio"
$ws.Range("H312").Value = "Americana (o)"
$ws.Range("I312").Value = "Segunda"
$ws.Range("J312").Value = 1340
$ws.Range("K312").Value = 6500
$ws.Range("L312").Value = 7000
$ws.Range("M312").Value = 6750
$ws.Range("N312").Value = "`$/docena de matas"
$ws.Range("O312").Value = "Provincia del Elquí"
$ws.Range("P312").Value = 1125
$ws.Range("Q312").Value = 6
$ws.Range("R312").Value = "Hortaliza"
